# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# The "municipio-nombre" column (column C) metadata rows were re-curated:
#   C2: iaest-measure:municipio-nombre  -> sdmx-dimension:refArea
#   C3: medida                         -> dim
#   C4: xsd:int                        -> URI-Municipio

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("C3").Value = "dim"
$ws.Range("C4").Value = "URI-Municipio"
